# Auto-generated by analysis script - Titan_Profits market data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 361.26666
$ws.Range("I53").Value = 262.9091
$ws.Range("J53").Value = 631.75
$ws.Range("K53").Value = 262.9091
$ws.Range("L53").Value = 631.75
$ws.Range("M53").Value = 374.0909
$ws.Range("N53").Value = -1905.75
$ws.Range("H113").Value = 2130.8333
$ws.Range("I113").Value = 2130.8333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2130.8333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1123.1667
$ws.Range("N113").ClearContents()
$ws.Range("H115").Value = 859.2308
$ws.Range("I115").Value = 680.8333
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 2042.4999
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -475.4999
$ws.Range("N115").Value = -12134
$ws.Range("H116").Value = 10646453
$ws.Range("I116").Value = 15377088
$ws.Range("K116").Value = 15377088
$ws.Range("M116").Value = -15373646
$ws.Range("H132").Value = 490356.88
$ws.Range("I132").Value = 578667.4399999999
$ws.Range("J132").Value = 26726.5
$ws.Range("K132").Value = 1736002.32
$ws.Range("L132").Value = 80179.5
$ws.Range("M132").Value = -1733472.32
$ws.Range("N132").Value = -85239.5
$ws.Range("H137").Value = 37038156
$ws.Range("I137").Value = 58824264
$ws.Range("K137").Value = 176472792
$ws.Range("M137").Value = -176470242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 30457.705
$ws.Range("I2").Value = 44568.22
$ws.Range("J2").Value = 953.9091
$ws.Range("K2").Value = 44568.22
$ws.Range("L2").Value = 953.9091
$ws.Range("M2").Value = -44455.22
$ws.Range("N2").Value = -1179.9091
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H63").Value = 9463.3125
$ws.Range("I63").Value = 9594.200000000001
$ws.Range("J63").Value = 7500
$ws.Range("K63").Value = 9594.200000000001
$ws.Range("L63").Value = 7500
$ws.Range("M63").Value = -8908.200000000001
$ws.Range("N63").Value = -8872
$ws.Range("H66").Value = 9463.3125
$ws.Range("I66").Value = 9594.200000000001
$ws.Range("J66").Value = 7500
$ws.Range("K66").Value = 47971
$ws.Range("L66").Value = 37500
$ws.Range("M66").Value = -44539
$ws.Range("N66").Value = -44364
$ws.Range("H116").Value = 30457.705
$ws.Range("I116").Value = 44568.22
$ws.Range("J116").Value = 953.9091
$ws.Range("K116").Value = 44568.22
$ws.Range("L116").Value = 953.9091
$ws.Range("M116").Value = -42274.22
$ws.Range("N116").Value = -5541.9091
$ws.Range("H122").Value = 3884
$ws.Range("I122").Value = 1756
$ws.Range("J122").Value = 6366.6665
$ws.Range("K122").Value = 5268
$ws.Range("L122").Value = 19099.9995
$ws.Range("M122").Value = -2818
$ws.Range("N122").Value = -23999.9995
$ws.Range("H132").Value = 2203.9524
$ws.Range("I132").Value = 1742.3889
$ws.Range("K132").Value = 5227.1667
$ws.Range("M132").Value = -2697.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 30457.705
$ws.Range("I3").Value = 44568.22
$ws.Range("J3").Value = 953.9091
$ws.Range("K3").Value = 44568.22
$ws.Range("L3").Value = 953.9091
$ws.Range("M3").Value = -44454.22
$ws.Range("N3").Value = -1181.9091
$ws.Range("H99").Value = 1301.96
$ws.Range("I99").Value = 1215.174
$ws.Range("J99").Value = 2300
$ws.Range("K99").Value = 1215.174
$ws.Range("L99").Value = 2300
$ws.Range("M99").Value = 282.826
$ws.Range("N99").Value = -5296

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1133.1666
$ws.Range("I16").Value = 1159.8
$ws.Range("K16").Value = 1159.8
$ws.Range("M16").Value = -872.8
$ws.Range("H31").Value = 3239.6333
$ws.Range("I31").Value = 1311.64
$ws.Range("J31").Value = 4616.7715
$ws.Range("K31").Value = 1311.64
$ws.Range("L31").Value = 4616.7715
$ws.Range("M31").Value = -1016.64
$ws.Range("N31").Value = -5206.7715
$ws.Range("H34").Value = 3239.6333
$ws.Range("I34").Value = 1311.64
$ws.Range("J34").Value = 4616.7715
$ws.Range("K34").Value = 1311.64
$ws.Range("L34").Value = 4616.7715
$ws.Range("M34").Value = -1109.64
$ws.Range("N34").Value = -5020.7715
$ws.Range("H58").Value = 3176.5789
$ws.Range("I58").Value = 1583.75
$ws.Range("J58").Value = 4335
$ws.Range("K58").Value = 1583.75
$ws.Range("L58").Value = 4335
$ws.Range("M58").Value = -1380.75
$ws.Range("N58").Value = -4741
$ws.Range("H76").Value = 2600
$ws.Range("I76").Value = 2600
$ws.Range("K76").Value = 2600
$ws.Range("M76").Value = -2285
$ws.Range("H79").Value = 2600
$ws.Range("I79").Value = 2600
$ws.Range("K79").Value = 2600
$ws.Range("M79").Value = -1508
$ws.Range("H105").Value = 1448.8
$ws.Range("I105").Value = 1312.5714
$ws.Range("J105").Value = 1766.6666
$ws.Range("K105").Value = 1312.5714
$ws.Range("L105").Value = 1766.6666
$ws.Range("M105").Value = 434.4286
$ws.Range("N105").Value = -5260.6666
$ws.Range("H113").Value = 1133.1666
$ws.Range("I113").Value = 1159.8
$ws.Range("K113").Value = 1159.8
$ws.Range("M113").Value = 1010.2
$ws.Range("H132").Value = 3036.0344
$ws.Range("I132").Value = 2617.6086
$ws.Range("K132").Value = 7852.825800000001
$ws.Range("M132").Value = -5322.825800000001
$ws.Range("H134").Value = 3976.889
$ws.Range("I134").Value = 1744.2
$ws.Range("K134").Value = 5232.6
$ws.Range("M134").Value = -2697.6
$ws.Range("H136").Value = 3176.5789
$ws.Range("I136").Value = 1583.75
$ws.Range("J136").Value = 4335
$ws.Range("K136").Value = 4751.25
$ws.Range("L136").Value = 13005
$ws.Range("M136").Value = -2201.25
$ws.Range("N136").Value = -18105

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76923260
$ws.Range("I2").Value = 15.428572
$ws.Range("K2").Value = 92.571432
$ws.Range("M2").Value = 20.428568
$ws.Range("H5").Value = 1614.4814
$ws.Range("J5").Value = 2683.077
$ws.Range("L5").Value = 8049.231000000001
$ws.Range("N5").Value = -8273.231
$ws.Range("H49").Value = 1402.5333
$ws.Range("J49").Value = 1504
$ws.Range("L49").Value = 4512
$ws.Range("N49").Value = -4824
$ws.Range("H70").Value = 3249.75
$ws.Range("I70").Value = 1999.5
$ws.Range("K70").Value = 5998.5
$ws.Range("M70").Value = -5683.5
$ws.Range("H73").Value = 3249.75
$ws.Range("I73").Value = 1999.5
$ws.Range("K73").Value = 5998.5
$ws.Range("M73").Value = -4906.5
$ws.Range("H76").Value = 3571
$ws.Range("I76").Value = 1713
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 5139
$ws.Range("L76").Value = 13500
$ws.Range("M76").Value = -4756
$ws.Range("N76").Value = -14266
$ws.Range("H79").Value = 3571
$ws.Range("I79").Value = 1713
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 5139
$ws.Range("L79").Value = 13500
$ws.Range("M79").Value = -3813
$ws.Range("N79").Value = -16152
$ws.Range("H131").Value = 5557335.5
$ws.Range("I131").Value = 229.5
$ws.Range("J131").Value = 6174791.5
$ws.Range("K131").Value = 688.5
$ws.Range("L131").Value = 18524374.5
$ws.Range("M131").Value = 4351.5
$ws.Range("N131").Value = -18534454.5
$ws.Range("H135").Value = 1614.4814
$ws.Range("J135").Value = 2683.077
$ws.Range("L135").Value = 24147.693
$ws.Range("N135").Value = -29217.693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2810.3333
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 3099
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 3099
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -7439
$ws.Range("H126").Value = 3020.5642
$ws.Range("I126").Value = 2080.1333
$ws.Range("J126").Value = 3608.3333
$ws.Range("K126").Value = 6240.3999
$ws.Range("L126").Value = 10824.9999
$ws.Range("M126").Value = -3770.3999
$ws.Range("N126").Value = -15764.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1259
$ws.Range("I46").Value = 1223.75
$ws.Range("K46").Value = 1223.75
$ws.Range("M46").Value = -1035.75
$ws.Range("H132").Value = 3580.9644
$ws.Range("I132").Value = 1934.6
$ws.Range("J132").Value = 5480.615
$ws.Range("K132").Value = 5803.799999999999
$ws.Range("L132").Value = 16441.845
$ws.Range("M132").Value = -3273.799999999999
$ws.Range("N132").Value = -21501.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4410.3213
$ws.Range("I132").Value = 4904.3335
$ws.Range("J132").Value = 2928.2856
$ws.Range("K132").Value = 14713.0005
$ws.Range("L132").Value = 8784.856800000001
$ws.Range("M132").Value = -12183.0005
$ws.Range("N132").Value = -13844.8568
$ws.Range("H136").Value = 2471.7917
$ws.Range("I136").Value = 1724.6154
$ws.Range("K136").Value = 5173.8462
$ws.Range("M136").Value = -2623.8462
